# Updated cryptos list with GitHub Actions
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row.
# Some Price values are plain decimals (e.g. "0.7009", "102.50"); a leading
# apostrophe is used (Excel's text quote-prefix) so they are stored as text
# with their exact original formatting (incl. trailing zeros) instead of
# being auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.218.20'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.854.42'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''0.7009'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').Value = '''237.51'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '''0.07934'
$ws.Range('E8').Value = '  +3.07%  '
$ws.Range('D9').Value = '''0.3016'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('D10').Value = '''23.54'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('D11').Value = '''0.08187'
$ws.Range('E11').Value = '  +2.07%  '
$ws.Range('D12').Value = '1.846.88'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').Value = '''5.180'
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '''0.7036'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').Value = '''89.46'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = '29.209.32'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '''5.810'
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('D18').Value = '''0.000007816'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = '''13.18'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '''236.02'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '2.088.95'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = '''1.002'
$ws.Range('D24').Value = '''7.488'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '''162.49'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').Value = '''8.850'
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('D27').Value = '''0.1411'
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').Value = '''18.06'
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').Value = '''1.909'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('D30').Value = '''1.404'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '''1.469'
$ws.Range('D32').Value = '''4.320'
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('D33').Value = '''4.011'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').Value = '''0.05148'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('D35').Value = '''1.163'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').Value = '''0.7098'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('D37').Value = '''0.9960'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('D38').Value = '''2.681'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').Value = '''0.01846'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('D40').Value = '''2.711'
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('D41').Value = '1.160.21'
$ws.Range('E41').Value = '  +5.26%  '
$ws.Range('D42').Value = '''0.9300'
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').Value = '''5.980'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').Value = '''0.4244'
$ws.Range('E44').Value = '  -0.93%  '
$ws.Range('D45').Value = '''70.08'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '''102.50'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '''0.5294'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('D49').Value = '''1.737'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').Value = '''9.133'
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').Value = '''6.953'
$ws.Range('E51').Value = '  -0.74%  '
